# "Restart experiment of 10-fold"
# Updates the accuracy figures in the 10-fold block (rows 10-19) of the
# single worksheet, moves the selection to C17, and normalises a couple
# of percentage number formats that lost their dedicated custom format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 10-fold results (rows 10-19) ---------------------------------------

# Fold 03 (row 10): both test accuracy and train accuracy change.
$ws.Range("B10").Value = 0.73469399999999996
$ws.Range("C10").NumberFormat = "0.0000%"
$ws.Range("C10").Value = 0.99587599999999998

# Fold 08 (row 11)
$ws.Range("B11").Value = 0.67241399999999996

# Fold 09 (row 12)
$ws.Range("B12").Value = 0.72092999999999996

# Fold 10 (row 13)
$ws.Range("B13").Value = 0.84210499999999999

# Fold 11 (row 14)
$ws.Range("B14").Value = 0.61818200000000001

# Fold 12 (row 15): train accuracy format changes from the old custom
# "0.0%" format to the "0.0000%" format used elsewhere in the sheet.
$ws.Range("C15").NumberFormat = "0.0000%"
$ws.Range("C15").Value = 0.99599199999999999

# Fold 13 (row 16)
$ws.Range("B16").Value = 0.88524599999999998

# Fold 14 (row 17)
$ws.Range("C17").Value = 0.96774199999999999

# Fold 15 (row 18): test accuracy format switches to the built-in "0%".
$ws.Range("B18").NumberFormat = "0%"
$ws.Range("B18").Value = 0.75

# Fold 16 (row 19): test accuracy format switches to "0.000%".
$ws.Range("B19").NumberFormat = "0.000%"
$ws.Range("B19").Value = 0.71831

# ---- selection -----------------------------------------------------------
$ws.Range("C17").Select() | Out-Null
